$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.205.51"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.144.27"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "534.64"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "139.03"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D8").Value = "3.139.76"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +4.48%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "0.108"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("D13").Value = "3.680.50"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "25.68"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "0.0000165"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "58.242.62"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "3.141.55"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +0.87%  "
$ws.Range("D20").Value = "12.74"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "8.18"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "360.98"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("D24").Value = "69.17"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").Value = "0.507"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").Value = "6.16"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "21.52"
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("D33").Value = "5.05"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").Value = "159.32"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").Value = "6.09"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "26.04"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("D40").Value = "0.0670"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "2.505.70"
$ws.Range("E41").Value = "  +8.62%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -4.33%  "
$ws.Range("D44").Value = "37.45"
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("D45").Value = "3.184.55"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0269"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "19.86"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "0.739"
$ws.Range("E51").Value = "  -4.27%  "
